{"js": "// Author's edit: the sentence describing how TCP/UDP file transfer was\n// implemented changed from\n//   \"...foram implementadas com a biblioteca \"socket\".\"\n// to\n//   \"...foi implementada atrav\u00e9s da API \"socket\".\"\n// and the stray \"_GoBack\" bookmark (left over from Word's last-edit-position\n// tracking, previously sitting mid-word inside the \"Enviar/Receber dados\n// (permanecer...\" bullet) moved to sit right before the re-typed \"socket\".\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Rewrite the core sentence fragment.\nconst sentenceHits = body.search(\"foram implementadas com a biblioteca\", { matchCase: true });\nsentenceHits.load(\"text\");\nawait context.sync();\n\nif (sentenceHits.items.length > 0) {\n  sentenceHits.items[0].insertText(\"foi implementada atrav\u00e9s da API\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Locate the paragraph still carrying the leftover \"_GoBack\" bookmark\n//    (two bullets share the same \"Enviar/Receber dados (pe...\" prefix, so\n//    disambiguate by checking which one actually owns the bookmark) and\n//    normalize it back to a single clean run with no bookmark.\nconst candidates = body.search(\"Enviar/Receber dados (pe\", { matchCase: true });\ncandidates.load(\"paragraphs\");\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < candidates.items.length; i++) {\n  const para = candidates.items[i].paragraphs.items[0];\n  const bookmarkResult = para.getRange().getBookmarks(true, true);\n  await context.sync();\n  if (bookmarkResult.value && bookmarkResult.value.indexOf(\"_GoBack\") !== -1) {\n    targetPara = para;\n    break;\n  }\n}\n\nif (targetPara) {\n  targetPara.load(\"text\");\n  await context.sync();\n  const mergedText = targetPara.text;\n\n  doc.deleteBookmark(\"_GoBack\");\n  // Force a real text delta (a same-text \"Replace\" is treated as a no-op and\n  // would leave the paragraph's runs split) so the paragraph re-collapses to\n  // a single run, then restore its original text.\n  targetPara.getRange().insertText(\"\\u0001\", \"Replace\");\n  await context.sync();\n  targetPara.getRange().insertText(mergedText, \"Replace\");\n  await context.sync();\n}\n\n// 3) Re-insert \"_GoBack\" right before the opening curly quote of \"socket\" in\n//    the freshly edited sentence.\nconst markHits = body.search(\"\u201csocket\u201d. Alguns detalhes\", { matchCase: true });\nmarkHits.load(\"text\");\nawait context.sync();\n\nif (markHits.items.length > 0) {\n  markHits.items[0].getRange(\"Start\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Author's edit: the sentence describing how TCP/UDP file transfer was\n# implemented changed from\n#   \"...foram implementadas com a biblioteca \"socket\".\"\n# to\n#   \"...foi implementada atrav\u00e9s da API \"socket\".\"\n# and the stray \"_GoBack\" bookmark (left over from Word's last-edit-position\n# tracking, previously sitting mid-word inside the \"Enviar/Receber dados\n# (permanecer...\" bullet) moved to sit right before the re-typed \"socket\".\n\n$d = $word.ActiveDocument\n\n# 1) Rewrite the core sentence fragment (find & replace).\n$range = $d.Content\n$null = $range.Find.Execute(\n    \"foram implementadas com a biblioteca\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"foi implementada atrav\u00e9s da API\",\n    2\n)\n\n# 2) Locate the paragraph still carrying the leftover \"_GoBack\" bookmark and\n#    normalize it back to a single clean run with no bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $bm = $d.Bookmarks(\"_GoBack\")\n    $bmPara = $bm.Range.Paragraphs(1)\n    $paraStart = $bmPara.Range.Start\n    $paraEnd = $bmPara.Range.End - 1   # exclude the trailing paragraph mark\n    $paraText = $d.Range($paraStart, $paraEnd).Text\n\n    $bm.Delete()\n\n    # A same-text assignment is treated as a no-op and would leave the\n    # paragraph's runs split across the old bookmark boundary, so force a\n    # real text delta first (collapse to a 1-char placeholder) and then\n    # restore the original text in one shot, landing it in a single run.\n    $tmpRange = $d.Range($paraStart, $paraEnd)\n    $tmpRange.Text = [char]1\n    $finalRange = $d.Range($paraStart, $paraStart + 1)\n    $finalRange.Text = $paraText\n}\n\n# 3) Re-insert \"_GoBack\" right before the opening curly quote of \"socket\" in\n#    the freshly edited sentence.\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\"\u201csocket\u201d. Alguns detalhes\")\nif ($found2) {\n    $markRange = $d.Range($range2.Start, $range2.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $markRange)\n}\n"}
